$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.699.57"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.724.57"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9972"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9974"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4930"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2610"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06224"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").Value = "1.728.34"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "15.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06993"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.500"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9982"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "26.502.00"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9972"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007195"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").Value = "1.943.33"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.440"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.562"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.102"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.399"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.744"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.917"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07991"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.667"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.607"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6260"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9372"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.001"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.426"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9969"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01512"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.575"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3856"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.909"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1159"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05383"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.805"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  -1.16%  "
